$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row containing "teste" (row 6), shifting "Pingo" (row 7) up to row 6
$ws.Rows.Item(6).Delete()
